# Scheduled runner update: refresh market-board pricing snapshot columns
# (currentAveragePrice / NQ / HQ / LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ)
# across the ALC, ARM, CRP, CUL, GSM, LTW and WVR leve-profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 337.6
$ws.Range("I6").Value = 334.9091
$ws.Range("J6").Value = 345
$ws.Range("K6").Value = 1004.7273
$ws.Range("L6").Value = 1035
$ws.Range("M6").Value = -892.7273
$ws.Range("N6").Value = -1259

$ws.Range("H12").Value = 553.2
$ws.Range("I12").Value = 581.44446
$ws.Range("K12").Value = 581.44446
$ws.Range("M12").Value = -411.44446

$ws.Range("H20").Value = 593.6667
$ws.Range("I20").Value = 593.6667
$ws.Range("K20").Value = 593.6667
$ws.Range("M20").Value = -363.6667

$ws.Range("H29").Value = 941.4
$ws.Range("I29").Value = 941.4
$ws.Range("K29").Value = 2824.2
$ws.Range("M29").Value = -2543.2

$ws.Range("H33").Value = 252.29411
$ws.Range("I33").Value = 263.4
$ws.Range("J33").Value = 169
$ws.Range("K33").Value = 263.4
$ws.Range("L33").Value = 169
$ws.Range("M33").Value = -34.39999999999998
$ws.Range("N33").Value = -627

$ws.Range("H35").Value = 593.6667
$ws.Range("I35").Value = 593.6667
$ws.Range("K35").Value = 593.6667
$ws.Range("M35").Value = -214.6667

$ws.Range("H70").Value = 3811.75
$ws.Range("I70").Value = 5118.8
$ws.Range("J70").Value = 1633.3334
$ws.Range("K70").Value = 15356.4
$ws.Range("L70").Value = 4900.0002
$ws.Range("M70").Value = -15086.4
$ws.Range("N70").Value = -5440.0002

$ws.Range("H73").Value = 3811.75
$ws.Range("I73").Value = 5118.8
$ws.Range("J73").Value = 1633.3334
$ws.Range("K73").Value = 15356.4
$ws.Range("L73").Value = 4900.0002
$ws.Range("M73").Value = -14420.4
$ws.Range("N73").Value = -6772.0002

$ws.Range("H74").Value = 5986
$ws.Range("I74").Value = 5986
$ws.Range("K74").Value = 5986
$ws.Range("M74").Value = -5050

$ws.Range("H77").Value = 5986
$ws.Range("I77").Value = 5986
$ws.Range("K77").Value = 29930
$ws.Range("M77").Value = -25250

$ws.Range("H106").Value = 7548.8125
$ws.Range("I106").Value = 6234.4287
$ws.Range("K106").Value = 6234.4287
$ws.Range("M106").Value = -5603.4287

$ws.Range("H112").Value = 3301.75
$ws.Range("J112").Value = 3302.5
$ws.Range("L112").Value = 9907.5
$ws.Range("N112").Value = -12123.5

$ws.Range("H137").Value = 2155.7273
$ws.Range("I137").Value = 2214.5
$ws.Range("J137").Value = 1999
$ws.Range("K137").Value = 6643.5
$ws.Range("L137").Value = 5997
$ws.Range("M137").Value = -4093.5
$ws.Range("N137").Value = -11097

$ws.Range("H141").Value = 10825.154
$ws.Range("I141").Value = 10329.75
$ws.Range("K141").Value = 30989.25
$ws.Range("M141").Value = -25809.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18032.076
$ws.Range("I32").Value = 18032.076
$ws.Range("K32").Value = 18032.076
$ws.Range("M32").Value = -17745.076

$ws.Range("H61").Value = 1497.7646
$ws.Range("I61").Value = 1497.2858
$ws.Range("K61").Value = 1497.2858
$ws.Range("M61").Value = -1285.2858

$ws.Range("H110").Value = 5419.2
$ws.Range("I110").Value = 6660
$ws.Range("J110").Value = 456
$ws.Range("K110").Value = 6660
$ws.Range("L110").Value = 456
$ws.Range("M110").Value = -4615
$ws.Range("N110").Value = -4546

$ws.Range("H122").Value = 1846.2667
$ws.Range("I122").Value = 1870.8572
$ws.Range("K122").Value = 5612.571599999999
$ws.Range("M122").Value = -3162.571599999999

$ws.Range("H135").Value = 75000
$ws.Range("J135").Value = 75000
$ws.Range("L135").Value = 75000
$ws.Range("N135").Value = -85140

$ws.Range("H136").Value = 1497.7646
$ws.Range("I136").Value = 1497.2858
$ws.Range("K136").Value = 4491.857400000001
$ws.Range("M136").Value = -1941.857400000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2123.2856
$ws.Range("I31").Value = 1660.5
$ws.Range("K31").Value = 1660.5
$ws.Range("M31").Value = -1365.5

$ws.Range("H34").Value = 2123.2856
$ws.Range("I34").Value = 1660.5
$ws.Range("K34").Value = 1660.5
$ws.Range("M34").Value = -1458.5

$ws.Range("H132").Value = 1924.3478
$ws.Range("I132").Value = 1949.4
$ws.Range("K132").Value = 5848.200000000001
$ws.Range("M132").Value = -3318.200000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H36").Value = 142.5
$ws.Range("I36").Value = 142.5
$ws.Range("K36").Value = 427.5
$ws.Range("M36").Value = -258.5

$ws.Range("H60").Value = 2255.8
$ws.Range("J60").Value = 2490
$ws.Range("L60").Value = 7470
$ws.Range("N60").Value = -7972

$ws.Range("H69").Value = 783.5714
$ws.Range("I69").Value = 827
$ws.Range("K69").Value = 2481
$ws.Range("M69").Value = -1670

$ws.Range("H72").Value = 783.5714
$ws.Range("I72").Value = 827
$ws.Range("K72").Value = 7443
$ws.Range("M72").Value = -3387

$ws.Range("H129").Value = 2333.9285
$ws.Range("J129").Value = 3301.2856
$ws.Range("L129").Value = 9903.856800000001
$ws.Range("N129").Value = -19903.8568

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 634
$ws.Range("I97").Value = 488.6
$ws.Range("J97").Value = 997.5
$ws.Range("K97").Value = 488.6
$ws.Range("L97").Value = 997.5
$ws.Range("M97").Value = 7.399999999999977
$ws.Range("N97").Value = -1989.5

$ws.Range("H132").Value = 2219.75
$ws.Range("I132").Value = 1968.1666
$ws.Range("K132").Value = 5904.4998
$ws.Range("M132").Value = -3374.4998

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H25").Value = 200
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").ClearContents()

$ws.Range("H122").Value = 4249.75
$ws.Range("I122").Value = 4199.8
$ws.Range("J122").Value = 4333
$ws.Range("K122").Value = 12599.4
$ws.Range("L122").Value = 12999
$ws.Range("M122").Value = -10149.4
$ws.Range("N122").Value = -17899

$ws.Range("H132").Value = 2913.4348
$ws.Range("I132").Value = 2863.5264
$ws.Range("K132").Value = 8590.5792
$ws.Range("M132").Value = -6060.5792

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 100
$ws.Range("I6").Value = 100
$ws.Range("K6").Value = 100
$ws.Range("M6").Value = 15

$ws.Range("H45").Value = 34528
$ws.Range("J45").Value = 37136.332
$ws.Range("L45").Value = 37136.332
$ws.Range("N45").Value = -38118.332
